$d = $word.ActiveDocument

# Position the insertion point at the end of the last paragraph's text,
# i.e. right before its trailing paragraph mark (so InsertXML appends new
# paragraphs after it instead of replacing the paragraph mark itself).
# Rebuild the range via Document.Range(...) (rather than mutating the
# Paragraph's own Range in place) so the collapsed point anchors correctly.
$lastPara = $d.Paragraphs.Last
$pr = $lastPara.Range
$r = $d.Range($pr.Start, $pr.End)
$r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1 ; shrink end by one char (the pilcrow)
$r.Collapse(0) | Out-Null      # wdCollapseEnd

$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = @"
<w:p $W><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Esta es la modificaci&#243;n </w:t></w:r></w:p><w:p $W><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p><w:p $W><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Fakfnfkscnc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
"@

$r.InsertXML($xml) | Out-Null
